$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the Antonio Luis Palomino Marimon (19935223) row 17 entirely - shifts
# everything below up by one row. This mirrors the real edit: one worker row
# was removed from the statement of account table.
$ws.Rows(17).Delete()

# --- Header block updates ---
# Cant. Trabajadores / Cant. Periodos counters
$ws.Range("E11").Value = 258128
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 4

# --- Data table updates (rows 16-20 after the row delete above) ---
# Row 16: Oscar Manuel Lopez Polo - salary updated
$ws.Range("G16").Value = 1423500

# Row 17 (was the old row 18: Antonio Luis... ) now becomes Abel Antonio's
# first period entry
$ws.Range("C17").Value = "72306068"
$ws.Range("D17").Value = "ABEL ANTONIO MUÃ?OZ CERDA"
$ws.Range("E17").Value = "2505"
$ws.Range("F17").Value = 43654
$ws.Range("G17").Value = 1423500

# Row 18 (was old row 19) - Abel Antonio period 2506
$ws.Range("E18").Value = "2506"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# Row 19 (was old row 20) - Abel Antonio period 2507
$ws.Range("E19").Value = "2507"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

# Row 20 (was old row 21, bottom-bordered row) - Abel Antonio period 2508
$ws.Range("E20").Value = "2508"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500
